# Refresh the "想去人数" (F column) counts across all four sheets to match
# the gh-pages data snapshot generated at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 329
$ws.Cells.Item(3, 6).Value = 224
$ws.Cells.Item(4, 6).Value = 545
$ws.Cells.Item(5, 6).Value = 1331
$ws.Cells.Item(6, 6).Value = 643
$ws.Cells.Item(7, 6).Value = 341
$ws.Cells.Item(8, 6).Value = 25
$ws.Cells.Item(9, 6).Value = 151
$ws.Cells.Item(10, 6).Value = 403
$ws.Cells.Item(11, 6).Value = 6125
$ws.Cells.Item(12, 6).Value = 109
$ws.Cells.Item(14, 6).Value = 1889
$ws.Cells.Item(15, 6).Value = 4587
$ws.Cells.Item(18, 6).Value = 306
$ws.Cells.Item(19, 6).Value = 5318
$ws.Cells.Item(20, 6).Value = 6963
$ws.Cells.Item(21, 6).Value = 145
$ws.Cells.Item(22, 6).Value = 1080
$ws.Cells.Item(23, 6).Value = 743
$ws.Cells.Item(24, 6).Value = 3950
$ws.Cells.Item(25, 6).Value = 537
$ws.Cells.Item(27, 6).Value = 221
$ws.Cells.Item(28, 6).Value = 143
$ws.Cells.Item(29, 6).Value = 1044
$ws.Cells.Item(30, 6).Value = 1480
$ws.Cells.Item(31, 6).Value = 541
$ws.Cells.Item(32, 6).Value = 666
$ws.Cells.Item(33, 6).Value = 1667
$ws.Cells.Item(34, 6).Value = 232
$ws.Cells.Item(35, 6).Value = 1852
$ws.Cells.Item(36, 6).Value = 40
$ws.Cells.Item(37, 6).Value = 1219
$ws.Cells.Item(38, 6).Value = 43
$ws.Cells.Item(40, 6).Value = 673
$ws.Cells.Item(42, 6).Value = 462
$ws.Cells.Item(43, 6).Value = 3627
$ws.Cells.Item(44, 6).Value = 155
$ws.Cells.Item(45, 6).Value = 335
$ws.Cells.Item(46, 6).Value = 437
$ws.Cells.Item(47, 6).Value = 19
$ws.Cells.Item(48, 6).Value = 86
$ws.Cells.Item(49, 6).Value = 3938

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 1256
$ws.Cells.Item(26, 6).Value = 52

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 4313

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 4313
$ws.Cells.Item(3, 6).Value = 329
$ws.Cells.Item(4, 6).Value = 1256
$ws.Cells.Item(7, 6).Value = 224
$ws.Cells.Item(8, 6).Value = 545
$ws.Cells.Item(10, 6).Value = 1331
$ws.Cells.Item(12, 6).Value = 643
$ws.Cells.Item(13, 6).Value = 341
$ws.Cells.Item(14, 6).Value = 25
$ws.Cells.Item(15, 6).Value = 151
$ws.Cells.Item(16, 6).Value = 403
$ws.Cells.Item(17, 6).Value = 109
$ws.Cells.Item(19, 6).Value = 1889
$ws.Cells.Item(20, 6).Value = 4588
$ws.Cells.Item(21, 6).Value = 5318
$ws.Cells.Item(22, 6).Value = 5318
$ws.Cells.Item(23, 6).Value = 145
$ws.Cells.Item(24, 6).Value = 1080
$ws.Cells.Item(25, 6).Value = 743
$ws.Cells.Item(26, 6).Value = 3950
$ws.Cells.Item(27, 6).Value = 537
$ws.Cells.Item(28, 6).Value = 221
$ws.Cells.Item(30, 6).Value = 143
$ws.Cells.Item(31, 6).Value = 1044
$ws.Cells.Item(32, 6).Value = 1480
$ws.Cells.Item(33, 6).Value = 541
$ws.Cells.Item(34, 6).Value = 666
$ws.Cells.Item(35, 6).Value = 1667
$ws.Cells.Item(36, 6).Value = 1852
$ws.Cells.Item(39, 6).Value = 673
$ws.Cells.Item(43, 6).Value = 3627
$ws.Cells.Item(45, 6).Value = 155
$ws.Cells.Item(46, 6).Value = 335
$ws.Cells.Item(47, 6).Value = 437
$ws.Cells.Item(48, 6).Value = 86
$ws.Cells.Item(50, 6).Value = 3938
